$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 6 - this shifts the existing row 6
#    ("Extension.valueReference" / "valueReference" slice) down to row 7,
#    preserving its values/format intact.
$ws.Rows.Item(6).Insert()

# 2. Give the freshly inserted (blank) row 6 the same cell style/format
#    used by the other data rows (border + top-aligned wrap text), by
#    copying formats down from row 5.
$ws.Range("A5:AJ5").Copy()
$ws.Range("A6:AJ6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Populate the new row 6 with the "Extension.value[x]" parent-element
#    values.
$ws.Range("A6").Value = "Extension.value[x]"
$ws.Range("E6").Value = "0"
$ws.Range("F6").Value = "1"
$ws.Range("J6").Value = "Reference`n"
$ws.Range("K6").Value = "Value of extension"
$ws.Range("L6").Value = "Value of extension - may be a resource or one of a constrained set of the data types (see Extensibility in the spec for list)."
$ws.Range("AA6").Value = "type:`$this}`n"
$ws.Range("AD6").Value = "closed"
$ws.Range("AE6").Value = "Extension.value[x]"
$ws.Range("AF6").Value = "0"
$ws.Range("AG6").Value = "1"
$ws.Range("AJ6").Value = "N/A"

# 4. Row 7 (the shifted former row 6) keeps all its old content except the
#    Path column, which now reads "Extension.value[x]" too (it became the
#    sliced child row under the new parent row).
$ws.Range("A7").Value = "Extension.value[x]"

# 5. Row 6/7 are hidden, like the rest of the detail rows - loading/saving
#    drops the row "hidden" flag unless re-asserted, so restate it for
#    every data row.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(5).Hidden = $true
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(7).Hidden = $true

# 6. Column A got narrower after the edit.
$ws.Columns.Item(1).ColumnWidth = 19.00390625

# 7. Re-point the defined name backing the autofilter database to the new
#    7-row extent.
foreach ($n in $wb.Names) {
    $n.RefersTo = "=Elements!`$A`$1:`$AJ`$7"
}

# 8. Rebuild the AutoFilter over the new A1:AJ7 extent with the same two
#    column filters it had before (col G <> " ", col AA blank).
$ws.AutoFilterMode = $false
$rng = $ws.Range("A1:AJ7")
[void]$rng.AutoFilter(7, "<>" + " ", 1)
[void]$rng.AutoFilter(27, @(""), 7)

# 9. Extend the conditional formatting that highlights the data rows so it
#    covers the newly added row.
$cf = $ws.Range("A2:AI5")
foreach ($fc in $cf.FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("A2:AI6"))
}

Write-Host "edit complete"
